$d = $word.ActiveDocument

# The paragraph currently contains the bold text split across two runs
# ("DOCX, DOC, PDF, HTML, XPS, R" + "TF and TXT") with a "_GoBack" bookmark
# sitting between them. Word's Find/Replace merges the text into a single
# contiguous run (preserving the bold formatting) and drops the now
# meaningless bookmark markers that separated the two runs.
$result = $d.Content.Find.Execute(
    "DOCX, DOC, PDF, HTML, XPS, RTF and TXT",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "DOCX, DOC, PDF, HTML, XPS, RTF and TXT",
    2)

if (-not $result) {
    throw "Could not find the target text to normalize."
}
